$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.495.23'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.105.95'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.27'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5242'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4558'
$ws.Range("E8").Value = '  +4.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '56.26'
$ws.Range("E9").Value = '  +7.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09011'
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.178'
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.31'
$ws.Range("E12").Value = '  -1.64%  '
$ws.Range("D13").Value = '2.100.44'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.853'
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.107'
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001182'
$ws.Range("E16").Value = '  +4.93%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '97.30'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.009'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06654'
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.19'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.308'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '30.565.12'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.365'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").Value = '2.347.74'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.39'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.46'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.525'
$ws.Range("E29").Value = '  -4.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.72'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.227'
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1071'
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.654'
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.360'
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.960'
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.33'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.900'
$ws.Range("E37").Value = '  +7.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02590'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2323'
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.67'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6889'
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.249'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.337'
$ws.Range("E44").Value = '  +5.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6439'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.06'
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.251'
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000349'
$ws.Range("E49").Value = '  +17.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.211'
$ws.Range("E50").Value = '  -2.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '82.98'
$ws.Range("E51").Value = '  +0.60%  '
